# Add an explicit "page break before: off" setting to every paragraph in
# the document body, as well as to the paragraph-level heading/title
# styles (Heading 1-6, Title, Subtitle) defined in the style sheet.
#
# This mirrors the OOXML diff, which inserts
#   <w:pageBreakBefore w:val="0"/>
# as the first child of every <w:pPr> in word/document.xml, and as a
# child of <w:pPr> (after keepNext/keepLines) for the non-default
# paragraph styles in word/styles.xml.

$d = $word.ActiveDocument

# --- Body paragraphs -------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = 0
}

# --- Paragraph styles (Heading 1-6, Title, Subtitle, ...) ------------
foreach ($s in $d.Styles) {
    # Type 1 == wdStyleTypeParagraph. Skip the built-in "Normal" style,
    # which has no explicit pPr in the style sheet and is left untouched.
    if ($s.Type -eq 1 -and $s.NameLocal -ne "Normal") {
        $s.ParagraphFormat.PageBreakBefore = 0
    }
}
